$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview sheet: row 3 is the eddf713a file - status changes in zh-cn (B3) and de-de (C3) columns
$wsOverview.Range("B3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("C3").Value = "Handed back: in sync with en-US"

# zh-cn sheet: row 3 is the eddf713a file - Status (C3) changes, Latest Handback DateTime (H3) updates
$wsZhCn.Range("C3").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("H3").Value = "2016-03-19 18:37:51"

# de-de sheet: row 3 is the eddf713a file - Status (C3) changes, Latest Handback DateTime (H3) updates
$wsDeDe.Range("C3").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("H3").Value = "2016-03-19 18:37:56"
